# Endpoints_and_testcases.xlsx edit:
#  - remove the "Body Type" column (old column C, with "application/json" values)
#  - this shifts old D(Purpose)->C, old E(Example/hyperlink)->D, old F(JSON Example)->E
#  - apply wrap-text formatting to the (now) Purpose column and to the JSON Example
#    column/header, matching the new narrower column widths
#  - re-point the hyperlinks (previously in column E) to their new column D location
#  - adjust row heights for the rows whose wrapped text now needs more vertical space
#  - update column widths and the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture existing hyperlink info (address/url/display) before touching columns ---
$hyperlinkInfo = @()
foreach ($hl in $ws.Hyperlinks) {
    $hyperlinkInfo += ,@($hl.Range.Address(), $hl.Address)
}

# --- remove the old "Body Type" column (C) ---
$ws.Columns("C").Delete()

# --- rebuild the hyperlinks at their new (shifted one column left) location ---
$ws.Hyperlinks.Delete()
foreach ($pair in $hyperlinkInfo) {
    $oldAddr = $pair[0]
    $url = $pair[1]
    $oldRange = $ws.Range($oldAddr)
    $newRange = $oldRange.Offset(0, -1)
    $ws.Hyperlinks.Add($newRange, $url) | Out-Null
}

# --- wrap text: header row (bold) cells C1 and E1 ---
$ws.Range("C1").WrapText = $true
$ws.Range("E1").WrapText = $true

# --- wrap text: Purpose column body (C2:C11) ---
$ws.Range("C2:C11").WrapText = $true

# --- wrap text: JSON Example column body (only populated rows) ---
$ws.Range("E5").WrapText = $true
$ws.Range("E8").WrapText = $true
$ws.Range("E11").WrapText = $true

# --- column widths ---
$ws.Columns("C").ColumnWidth = 29
$ws.Columns("E").ColumnWidth = 21.6

# --- row heights (auto-fit equivalent for the now-wrapped, narrower columns) ---
$ws.Rows(2).RowHeight = 28.8
$ws.Rows(3).RowHeight = 28.8
$ws.Rows(4).RowHeight = 28.8
$ws.Rows(5).RowHeight = 57.6
$ws.Rows(6).RowHeight = 28.8
$ws.Rows(7).RowHeight = 28.8
$ws.Rows(8).RowHeight = 72
$ws.Rows(9).RowHeight = 28.8
$ws.Rows(10).RowHeight = 28.8
$ws.Rows(11).RowHeight = 57.6

# --- selection / view ---
$ws.Range("D8").Select()

Write-Output "done"
